# Add "Day 7: Handy Haversacks" data to the "2020" personal-stats sheet.
# This fills in row 11 (previously a blank "Day 7: " placeholder row) with
# real puzzle times, and nudges the cached C/E/F numbers for the already
# completed days (rows 5-9) which get re-based slightly in this commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2020")

# --- Row 5 : Day 1 -------------------------------------------------
$ws.Range("C5").Value = 0.0048495370370370368
$ws.Range("E5").Value = 0.0053125000000000004
$ws.Range("F5").Value = 0.0053125000000000004

# --- Row 6 : Day 2 -------------------------------------------------
$ws.Range("C6").Value = 0.0029861111111111113
$ws.Range("E6").Value = 0.0048726851851851856
$ws.Range("F6").Value = 0.0042824074074074075

# --- Row 7 : Day 3 -------------------------------------------------
$ws.Range("C7").Value = 0.0063888888888888884
$ws.Range("E7").Value = 0.008726851851851852
$ws.Range("F7").Value = 0.004386574074074074

# --- Row 8 : Day 4 -------------------------------------------------
$ws.Range("C8").Value = 0.0091435185185185178
$ws.Range("E8").Value = 0.021736111111111112
$ws.Range("F8").Value = 0.012777777777777777

# --- Row 9 : Day 5 -------------------------------------------------
$ws.Range("C9").Value = 0.0095833333333333343
$ws.Range("E9").Value = 0.011932870370370371
$ws.Range("F9").Value = 0.0070717592592592594

# --- Row 10 : Day 6 --------------------------------------------------
$ws.Range("C10").Value = 0.0088657407407407417
$ws.Range("E10").Value = 0.010150462962962964
$ws.Range("F10").Value = 0.0042592592592592595

# --- Row 11 : Day 7 (new) -------------------------------------------
$ws.Range("B11").Value = "Day 7: Handy Haversacks"
$ws.Range("C11").Value = 0.0098726851851851857
$ws.Range("E11").Value = 0.020682870370370372
$ws.Range("F11").Value = 0.012199074074074072
$ws.Range("H11").Value = "4th"

# Move the active selection to match the author's final cursor position.
$ws.Range("F12").Select()
